$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 400
$ws.Range("J7").Value = 400
$ws.Range("L7").Value = 400
$ws.Range("N7").Value = -624

$ws.Range("H8").Value = 333333380
$ws.Range("I8").Value = 333333380
$ws.Range("K8").Value = 1000000140
$ws.Range("M8").Value = -1000000001

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H14").Value = 400
$ws.Range("J14").Value = 400
$ws.Range("L14").Value = 400
$ws.Range("N14").Value = -782

$ws.Range("H16").Value = 2966.6667
$ws.Range("I16").Value = 2966.6667
$ws.Range("K16").Value = 2966.6667
$ws.Range("M16").Value = -2736.6667

$ws.Range("H18").Value = 759
$ws.Range("I18").Value = 761.25
$ws.Range("J18").Value = 750
$ws.Range("K18").Value = 761.25
$ws.Range("L18").Value = 750
$ws.Range("M18").Value = -477.25
$ws.Range("N18").Value = -1318

$ws.Range("H21").Value = 8295.666999999999
$ws.Range("J21").Value = 8759.5
$ws.Range("L21").Value = 8759.5
$ws.Range("N21").Value = -9695.5

$ws.Range("H23").Value = 8295.666999999999
$ws.Range("J23").Value = 8759.5
$ws.Range("L23").Value = 8759.5
$ws.Range("N23").Value = -9227.5

$ws.Range("H76").Value = 8521.519
$ws.Range("I76").Value = 10317.5625
$ws.Range("J76").Value = 5909.091
$ws.Range("K76").Value = 10317.5625
$ws.Range("L76").Value = 5909.091
$ws.Range("M76").Value = -10002.5625
$ws.Range("N76").Value = -6539.091

$ws.Range("H79").Value = 8521.519
$ws.Range("I79").Value = 10317.5625
$ws.Range("J79").Value = 5909.091
$ws.Range("K79").Value = 10317.5625
$ws.Range("L79").Value = 5909.091
$ws.Range("M79").Value = -9225.5625
$ws.Range("N79").Value = -8093.091

$ws.Range("H113").Value = 3129.0588
$ws.Range("I113").Value = 3400.375
$ws.Range("J113").Value = 3045.577
$ws.Range("K113").Value = 3400.375
$ws.Range("L113").Value = 3045.577
$ws.Range("M113").Value = -146.375
$ws.Range("N113").Value = -9553.577000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 500
$ws.Range("J8").Value = 500
$ws.Range("L8").Value = 500
$ws.Range("N8").Value = -788

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").ClearContents()
$ws.Range("N11").Value = 0

$ws.Range("H13").Value = 20000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 20000
$ws.Range("K13").Value = 0
$ws.Range("L13").ClearContents()
$ws.Range("M13").Value = 20000
$ws.Range("N13").Value = -20288

$ws.Range("H61").Value = 3099.682
$ws.Range("I61").Value = 2091.5
$ws.Range("J61").Value = 4309.5
$ws.Range("K61").Value = 2091.5
$ws.Range("L61").Value = 4309.5
$ws.Range("M61").Value = -1879.5
$ws.Range("N61").Value = -4733.5

$ws.Range("H136").Value = 3099.682
$ws.Range("I136").Value = 2091.5
$ws.Range("J136").Value = 4309.5
$ws.Range("K136").Value = 6274.5
$ws.Range("L136").Value = 12928.5
$ws.Range("M136").Value = -3724.5
$ws.Range("N136").Value = -18028.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("N14").Value = 0

$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("N15").Value = 0

$ws.Range("H17").Value = 725.75
$ws.Range("I17").Value = 554
$ws.Range("J17").Value = 897.5
$ws.Range("K17").Value = 554
$ws.Range("L17").Value = 897.5
$ws.Range("M17").Value = -382
$ws.Range("N17").Value = -1241.5

$ws.Range("H22").Value = 385.7143
$ws.Range("I22").Value = 260
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 260
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -87
$ws.Range("N22").Value = -1046

$ws.Range("H94").Value = 916.375
$ws.Range("I94").Value = 568.23334
$ws.Range("J94").Value = 1960.8
$ws.Range("K94").Value = 568.23334
$ws.Range("L94").Value = 1960.8
$ws.Range("M94").Value = -117.23334
$ws.Range("N94").Value = -2862.8

$ws.Range("H134").Value = 29899.15
$ws.Range("I134").Value = 39126.285
$ws.Range("J134").Value = 8369.166999999999
$ws.Range("K134").Value = 117378.855
$ws.Range("L134").Value = 25107.501
$ws.Range("M134").Value = -114843.855
$ws.Range("N134").Value = -30177.501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3178.4
$ws.Range("I31").Value = 2490.077
$ws.Range("J31").Value = 3637.282
$ws.Range("K31").Value = 2490.077
$ws.Range("L31").Value = 3637.282
$ws.Range("M31").Value = -2195.077
$ws.Range("N31").Value = -4227.282

$ws.Range("H34").Value = 3178.4
$ws.Range("I34").Value = 2490.077
$ws.Range("J34").Value = 3637.282
$ws.Range("K34").Value = 2490.077
$ws.Range("L34").Value = 3637.282
$ws.Range("M34").Value = -2288.077
$ws.Range("N34").Value = -4041.282

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 804.8
$ws.Range("I4").Value = 110
$ws.Range("J4").Value = 1268
$ws.Range("K4").Value = 330
$ws.Range("L4").Value = 3804
$ws.Range("M4").Value = -218
$ws.Range("N4").Value = -4028

$ws.Range("H9").Value = 39376124
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 39376124
$ws.Range("K9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("M9").Value = 118128372
$ws.Range("N9").Value = -118128820

$ws.Range("H15").Value = 257.2
$ws.Range("I15").Value = 60.4
$ws.Range("J15").Value = 454
$ws.Range("K15").Value = 181.2
$ws.Range("L15").Value = 1362
$ws.Range("M15").Value = -41.19999999999999
$ws.Range("N15").Value = -1642

$ws.Range("H40").Value = 249.09091
$ws.Range("I40").Value = 174
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 696
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -627
$ws.Range("N40").Value = -4138

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 479.25
$ws.Range("I9").Value = 479.25
$ws.Range("K9").Value = 479.25
$ws.Range("M9").Value = -309.25

$ws.Range("H12").Value = 980
$ws.Range("I12").Value = 980
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 980
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -840

$ws.Range("H18").Value = 6504938
$ws.Range("I18").Value = 8669918
$ws.Range("J18").Value = 10000
$ws.Range("K18").Value = 8669918
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = -8669625
$ws.Range("N18").Value = -10586

$ws.Range("H80").Value = 4956.6294
$ws.Range("I80").Value = 5305.6523
$ws.Range("K80").Value = 5305.6523
$ws.Range("M80").Value = -4307.6523

$ws.Range("H83").Value = 4956.6294
$ws.Range("I83").Value = 5305.6523
$ws.Range("K83").Value = 26528.2615
$ws.Range("M83").Value = -21536.2615

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 2374.25
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2374.25
$ws.Range("K13").Value = 0
$ws.Range("L13").ClearContents()
$ws.Range("M13").Value = 2374.25
$ws.Range("N13").Value = -2654.25

$ws.Range("H20").Value = 8333.333000000001
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 8333.333000000001
$ws.Range("K20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").Value = 8333.333000000001
$ws.Range("N20").Value = -8785.333000000001

$ws.Range("H132").Value = 5988.2085
$ws.Range("I132").Value = 1631.1562
$ws.Range("J132").Value = 14702.3125
$ws.Range("K132").Value = 4893.4686
$ws.Range("L132").Value = 44106.9375
$ws.Range("M132").Value = -2363.4686
$ws.Range("N132").Value = -49166.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 28000
$ws.Range("J58").Value = 28000
$ws.Range("L58").Value = 28000
$ws.Range("N58").Value = -28616

$ws.Range("H81").Value = 4713.3335
$ws.Range("J81").Value = 5375
$ws.Range("L81").Value = 10750
$ws.Range("N81").Value = -12872

$ws.Range("H84").Value = 4713.3335
$ws.Range("J84").Value = 5375
$ws.Range("L84").Value = 53750
$ws.Range("N84").Value = -64358

